# Regenerate merged AHB files
# - Rename the "_old"/"_new" header-name suffixes to "_FV2210"/"_FV2304"
# - Turn the header range A1:U70 into a real Excel Table ("Table1")
# - Freeze the header row (row 1) in the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header labels (columns A-J: *_old -> *_FV2210, columns L-U: *_new -> *_FV2304) ---
# Use ordered arrays (not a hashtable) so the write order is explicit/deterministic.
$renames = @(
    ,@("A1", "Segmentname_FV2210")
    ,@("B1", "Segmentgruppe_FV2210")
    ,@("C1", "Segment_FV2210")
    ,@("D1", "Datenelement_FV2210")
    ,@("E1", "Segment ID_FV2210")
    ,@("F1", "Code_FV2210")
    ,@("G1", "Qualifier_FV2210")
    ,@("H1", "Beschreibung_FV2210")
    ,@("I1", "Bedingungsausdruck_FV2210")
    ,@("J1", "Bedingung_FV2210")
    # K1 ("diff") is unchanged
    ,@("L1", "Segmentname_FV2304")
    ,@("M1", "Segmentgruppe_FV2304")
    ,@("N1", "Segment_FV2304")
    ,@("O1", "Datenelement_FV2304")
    ,@("P1", "Segment ID_FV2304")
    ,@("Q1", "Code_FV2304")
    ,@("R1", "Qualifier_FV2304")
    ,@("S1", "Beschreibung_FV2304")
    ,@("T1", "Bedingungsausdruck_FV2304")
    ,@("U1", "Bedingung_FV2304")
)
foreach ($pair in $renames) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# --- 2. Convert the A1:U70 range into an Excel Table ---
$tableRange = $ws.Range("A1:U70")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# --- 3. Freeze panes at row 2 (keep header row 1 visible) ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
